$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newBullets = "RGBW light engine with app + remote control (solid colors & effects).;`nTwinkle / meteor modes and optional music-sync animations.;`nHundreds of fiber strands installed for dense, even star coverage.;`nHidden wiring and fused power with OEM-style protection.;`nInstall typically completed in 2 days for sedans (vehicle-dependent).;`nWarranty on workmanship. "

$ws.Range("G6").Value = $newBullets

$ws.Range("C6").Select()
